$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1:B2").NumberFormat = "@"

$ws.Range("A1").Value = "81405441600015"
$ws.Range("B1").Value = "500,00 €"
$ws.Range("A2").Value = "81405440000000.00"
$ws.Range("B2").Value = "500.00"
